# Update "想去人数" (attendance / interest count) figures in column F
# for the two sheets that hold the full event data: "展览" and "全部类型".
# (The other two sheets, "演出" and "本地生活", only contain header rows
# and no data, so they are left untouched.)

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 7055
    6  = 4080
    11 = 94
    12 = 57
    13 = 47
    15 = 609
    16 = 82
    18 = 130
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
